# Refresh crypto price/volume data (D: Price, E: Volume(1h)) for rows 2-51.
# Cells in column D whose new value parses as a plain number must be forced
# to remain Text (matching the sheet's original inlineStr string cells) by
# temporarily applying a Text number format, then clearing formats again so
# no new cell style is introduced (preserves the original style index 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.738.14"
$ws.Range("D3").Value = "2.100.73"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.76"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.39"
$ws.Range("D7").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.74"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").Value = "2.412.64"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.06"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.810"
$ws.Range("D15").ClearFormats()
$ws.Range("D17").Value = "2.107.77"
$ws.Range("D18").Value = "38.756.53"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.76"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.10"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("D25").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "172.08"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.138"
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.34"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.52"
$ws.Range("D31").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.03"
$ws.Range("D35").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.55"
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.13"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.89"
$ws.Range("D41").ClearFormats()
$ws.Range("D43").Value = "1.530.28"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.80"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.13"
$ws.Range("D48").ClearFormats()
$ws.Range("D51").Value = "2.299.28"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +1.92%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("E11").Value = "  -1.19%  "
$ws.Range("E12").Value = "  +6.37%  "
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("E15").Value = "  +3.64%  "
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  +3.20%  "
$ws.Range("E29").Value = "  +2.76%  "
$ws.Range("E30").Value = "  +1.11%  "
$ws.Range("E31").Value = "  +8.69%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("E35").Value = "  +7.59%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  -2.60%  "
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("E44").Value = "  +6.67%  "
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("E51").Value = "  +0.53%  "
